# Regenerate merged AHB files
# - Rename header row 1 columns: "..._old" -> "..._FV2310", "..._new" -> "..._FV2404"
# - Turn the used range A1:U64 into an Excel Table ("Table1")
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header cells (columns A-J: _old -> _FV2310, L-U: _new -> _FV2404; K "diff" untouched) ---
$headerRenames = @{
    1  = "Segmentname_FV2310"
    2  = "Segmentgruppe_FV2310"
    3  = "Segment_FV2310"
    4  = "Datenelement_FV2310"
    5  = "Segment ID_FV2310"
    6  = "Code_FV2310"
    7  = "Qualifier_FV2310"
    8  = "Beschreibung_FV2310"
    9  = "Bedingungsausdruck_FV2310"
    10 = "Bedingung_FV2310"
    12 = "Segmentname_FV2404"
    13 = "Segmentgruppe_FV2404"
    14 = "Segment_FV2404"
    15 = "Datenelement_FV2404"
    16 = "Segment ID_FV2404"
    17 = "Code_FV2404"
    18 = "Qualifier_FV2404"
    19 = "Beschreibung_FV2404"
    20 = "Bedingungsausdruck_FV2404"
    21 = "Bedingung_FV2404"
}

foreach ($col in $headerRenames.Keys) {
    $ws.Cells.Item(1, $col).Value2 = $headerRenames[$col]
}

# --- 2) Convert the used range into a native Excel Table (ListObject) ---
$usedRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $usedRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"

# --- 3) Freeze panes at the header row (row 1) ---
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
